$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AE2").Value = 11.77990626116324
$ws.Range("AG2").Value = 9.402482728738628
$ws.Range("AE3").Value = 14.09549045834902
$ws.Range("AG3").Value = 9.708036622612987
$ws.Range("AE4").Value = 12.05507120829597
$ws.Range("AG4").Value = 7.81534108263544
$ws.Range("AE5").Value = 15.14395095318842
$ws.Range("AG5").Value = 10.87695782895565
$ws.Range("AE6").Value = 11.41663563776412
$ws.Range("AG6").Value = 9.117486087690445
$ws.Range("AE7").Value = 11.5961704773125
$ws.Range("AG7").Value = 6.765597240361365
$ws.Range("AE8").Value = 11.52701059773078
$ws.Range("AG8").Value = 5.417187475597562
$ws.Range("AE9").Value = 11.96681525992708
$ws.Range("AG9").Value = 5.075922352128496
$ws.Range("AE10").Value = 9.860478640744111
$ws.Range("AG10").Value = 5.378310113809826
$ws.Range("AE11").Value = 10.60976684966777
$ws.Range("AG11").Value = 4.660998636527554
$ws.Range("AE12").Value = 13.2386734764729
$ws.Range("AG12").Value = 9.432043845850497
$ws.Range("AE13").Value = 11.88895917226999
$ws.Range("AG13").Value = 10.05325969952986
$ws.Range("AE14").Value = 10.16256479919443
$ws.Range("AG14").Value = 4.431569518760055
$ws.Range("AE15").Value = 14.79637055232094
$ws.Range("AG15").Value = 12.61661663544376
$ws.Range("AE16").Value = 11.76908737648015
$ws.Range("AG16").Value = 8.056823076801912
$ws.Range("AE17").Value = 14.3740775733138
$ws.Range("AG17").Value = 10.65325492887831
$ws.Range("AE18").Value = 11.25392851252065
$ws.Range("AG18").Value = 6.050565869735037
$ws.Range("AE19").Value = 9.881465367855055
$ws.Range("AG19").Value = 5.087764911852834
$ws.Range("AE20").Value = 11.99716097123871
$ws.Range("AG20").Value = 7.179193988801213
$ws.Range("AE21").Value = 10.64084482047585
$ws.Range("AG21").Value = 3.040592930114375
$ws.Range("AE22").Value = 12.4497784749184
$ws.Range("AG22").Value = 7.626863474768419
$ws.Range("AE23").Value = 13.03096986013981
$ws.Range("AG23").Value = 8.293086482866025
$ws.Range("AE24").Value = 13.35229821660199
$ws.Range("AG24").Value = 7.650730966519666
$ws.Range("AE25").Value = 14.09205416163184
$ws.Range("AG25").Value = 9.615789593193673
$ws.Range("AE26").Value = 13.1559824118403
$ws.Range("AG26").Value = 9.494878669991595
$ws.Range("AE27").Value = 11.40694298563118
$ws.Range("AG27").Value = 7.677335405662255
$ws.Range("AE28").Value = 14.90906696606726
$ws.Range("AG28").Value = 9.456186442353856
$ws.Range("AE29").Value = 10.6625181294375
$ws.Range("AG29").Value = 7.252800152167314
$ws.Range("AE30").Value = 12.26650352013645
$ws.Range("AG30").Value = 6.088463038734886
$ws.Range("AE31").Value = 12.54315067882748
$ws.Range("AG31").Value = 5.707182843100437
$ws.Range("AE32").Value = 14.04173456276912
$ws.Range("AG32").Value = 10.78667113646602
$ws.Range("AE33").Value = 10.61390266985545
$ws.Range("AG33").Value = 4.643561234844674
$ws.Range("AE34").Value = 13.34067530921292
$ws.Range("AG34").Value = 9.338360308913188
$ws.Range("AE35").Value = 10.68295423017317
$ws.Range("AG35").Value = 4.83028314597836
$ws.Range("AE36").Value = 11.11820274898126
$ws.Range("AG36").Value = 5.866352520342057
$ws.Range("AE37").Value = 12.19860116682559
$ws.Range("AG37").Value = 7.882414368586492
$ws.Range("AE38").Value = 13.2024270099755
$ws.Range("AG38").Value = 9.653347031077505
$ws.Range("AE39").Value = 12.9997840726113
$ws.Range("AG39").Value = 9.659659405851659
$ws.Range("AE40").Value = 15.06446348441912
$ws.Range("AG40").Value = 10.79900850881555
$ws.Range("AE41").Value = 13.48522943201806
$ws.Range("AG41").Value = 8.732100402036355
$ws.Range("AE42").Value = 12.67359967908071
$ws.Range("AG42").Value = 8.310303314519064
$ws.Range("AE43").Value = 11.90477475879558
$ws.Range("AG43").Value = 5.890685239368419
$ws.Range("AE44").Value = 13.79892655811657
$ws.Range("AG44").Value = 8.384466384050224
$ws.Range("AE45").Value = 12.93795781574384
$ws.Range("AG45").Value = 7.397599474703376
$ws.Range("AE46").Value = 14.14625678105321
$ws.Range("AG46").Value = 9.999666935118119
